# repull data, push all data, mean calculation
# Update column F (dSF) values on Sheet1 for the rows whose source data
# was repulled, per the upstream diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$updates = @{
    "F2"  = -7
    "F3"  = 3
    "F4"  = 1
    "F6"  = 3
    "F7"  = 6
    "F8"  = 5
    "F9"  = -8
    "F10" = -5
    "F13" = 1
    "F14" = 2
    "F15" = -3
    "F16" = -2
    "F17" = -3
    "F18" = -1
    "F19" = -7
    "F23" = -2
    "F25" = -2
    "F27" = 0
    "F28" = 1
    "F30" = -3
}

foreach ($cellRef in $updates.Keys) {
    $ws.Range($cellRef).Value = $updates[$cellRef]
}
